# Nexial "electron-showcase" workbook update:
#   [base] - [outputToCloud(resource)]: support the transferring of output
#   artifact to the cloud.
#
# This adds:
#   * outputToCloud(resource) to the "base" command category (column E).
#   * a brand new "text" command category (new column Y) containing
#     spellCheck(var,profile,text), and registers "text" in the "target"
#     category list (column A).
#
# The "#system" sheet backs a bunch of named ranges that feed dropdown /
# autocomplete validation lists elsewhere in the workbook, so the existing
# columns Y..AD (web, webalert, webcookie, ws, ws.async, xml) need to slide
# one column to the right (Z..AE) to make room for the new "text" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# Helper: shift every value of one column into the column immediately to
# its right, working bottom-most-affected-column first so we never clobber
# data before it has been copied.
# ---------------------------------------------------------------------------
function Shift-ColumnRight($fromCol, $toCol, $lastRow) {
    for ($r = 1; $r -le $lastRow; $r++) {
        $v = $ws.Range("$fromCol$r").Value2
        $ws.Range("$toCol$r").Value = $v
    }
}

Shift-ColumnRight "AD" "AE" 27
Shift-ColumnRight "AC" "AD" 8
Shift-ColumnRight "AB" "AC" 17
Shift-ColumnRight "AA" "AB" 8
Shift-ColumnRight "Z"  "AA" 8
Shift-ColumnRight "Y"  "Z"  129

# ---------------------------------------------------------------------------
# New "text" category: header + single command, replacing the vacated
# column Y.
# ---------------------------------------------------------------------------
$ws.Range("Y1").Value = "text"
$ws.Range("Y2").Value = "spellCheck(var,profile,text)"

# ---------------------------------------------------------------------------
# "target" category list (column A): insert "text" alphabetically between
# "step" and "web", pushing the remaining entries down by one row.
# ---------------------------------------------------------------------------
for ($r = 30; $r -ge 25; $r--) {
    $v = $ws.Range("A$r").Value2
    $ws.Range("A" + ($r + 1)).Value = $v
}
$ws.Range("A25").Value = "text"

# ---------------------------------------------------------------------------
# "base" category list (column E): insert outputToCloud(resource)
# alphabetically between "macro(file,sheet,name)" and
# "prependText(var,prependWith)", pushing the remaining entries down by one
# row.
# ---------------------------------------------------------------------------
for ($r = 38; $r -ge 22; $r--) {
    $v = $ws.Range("E$r").Value2
    $ws.Range("E" + ($r + 1)).Value = $v
}
$ws.Range("E22").Value = "outputToCloud(resource)"

# ---------------------------------------------------------------------------
# Update the named ranges that describe each category list so they keep
# pointing at the right (now one-row-longer / one-column-shifted) range.
# ---------------------------------------------------------------------------
$wb.Names.Item("base").RefersTo       = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo     = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo        = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo   = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo  = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo         = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo   = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo        = "='#system'!`$AE`$2:`$AE`$27"
$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")
